$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...the second child of Robert Pfitzner..." ->
#           "...the second child of Carl Robert Pfitzner..."
# Insert "Carl " right before "Robert Pfitzner, a violinist" as its own run.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Robert Pfitzner, a violinist", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s1 = $rng1.Start
$ins1 = $d.Range($s1, $s1)
$ins1.InsertBefore("Carl ")
# Force the newly typed text to live in its own run (rather than being
# silently re-merged with its neighbours) by toggling a character property.
$carlRange = $d.Range($s1, $s1 + "Carl ".Length)
$carlRange.Font.Bold = 1
$carlRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# Change 2: add a new dated entry "27 September 1869  Felipe Pedrell marries
# Carmen Domingo." right before the existing "28 September 1869" entry.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("28 September 1869", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s2 = $rng2.Start
$insPara = $d.Range($s2, $s2)
$insPara.InsertParagraphBefore()

$dateText = "27 September 1869"
$bodyText = "  Felipe Pedrell marries Carmen Domingo."

$dateRange = $d.Range($s2, $s2)
$dateRange.InsertAfter($dateText)
$dateEnd = $s2 + $dateText.Length
$boldRange = $d.Range($s2, $dateEnd)
$boldRange.Font.Bold = 1

$bodyRange = $d.Range($dateEnd, $dateEnd)
$bodyRange.InsertAfter($bodyText)
$bodyEnd = $dateEnd + $bodyText.Length
$plainRange = $d.Range($dateEnd, $bodyEnd)
$plainRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# Change 3: the closing date "12 March 2016" -> "12 June 2016"
# (keep the "12 " run intact; split "March 2016" into "June" + " 2016").
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Paul Scharfenberger", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterCredit = $d.Range($rng3.End, $d.Content.End)
$afterCredit.Find.Execute("March 2016", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s3 = $afterCredit.Start
$monthRange = $d.Range($s3, $s3 + 5)
$monthRange.Font.Bold = 1
$monthRange.Text = "June"
$juneEnd = $s3 + "June".Length
$juneRange = $d.Range($s3, $juneEnd)
$juneRange.Font.Bold = 0
